$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list on Fri Mar  8 17:24:58 UTC 2024 with GitHub Actions
# Each row: Coin (B), Link (C), Price (D), Volume(1h) (E)

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.858.05"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.36%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.911.29"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.71%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "479.47"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.60%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.00"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.27%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.620"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.45%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.997"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.14%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.723"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.71%  "

# Row 10
$ws.Range("E10").Value = "  +7.38%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000351"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +10.15%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.37"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.01%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.523.50"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.55%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.34"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.08%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.62"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.04%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.926.00"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.76%  "

# Row 17
$ws.Range("E17").Value = "  -0.29%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.77"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.68%  "

# Row 19
$ws.Range("E19").Value = "  -2.85%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.001.39"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.42%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "431.74"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.51%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.40"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.70%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.44"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.17%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.99"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.43%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.81"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +5.45%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.55"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.43%  "

# Row 27
$ws.Range("E27").Value = "  +2.33%  "

# Row 28
$ws.Range("E28").Value = "  +8.24%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.41"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +6.54%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "721.25"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.30%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.24"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.82%  "

# Row 32
$ws.Range("E32").Value = "  -4.26%  "

# Row 33
$ws.Range("E33").Value = "  +3.80%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0890"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +29.47%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "41.68"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.40%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.98"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.58%  "

# Row 37
$ws.Range("E37").Value = "  -6.18%  "

# Row 38
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.46"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.70%  "

# Row 39
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.15%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.88"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +8.80%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0468"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.96%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.01"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +11.20%  "

# Row 43
$ws.Range("E43").Value = "  +1.18%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.347"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.43%  "

# Row 45
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.140"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.09%  "

# Row 46
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.19%  "

# Row 47
$ws.Range("E47").Value = "  -0.25%  "

# Row 48
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.17"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.51%  "

# Row 49
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.23"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.76%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "145.24"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.10%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.84"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.32%  "
